# Applies the "cryptos list" price/volume refresh described by the commit:
#   "Updated cryptos list on Mon Jun 19 06:50:35 UTC 2023 with GitHub Actions"
#
# - Columns D (Price) and E (Volume 1h) are refreshed to new scraped values for
#   most rows (2-51). These are free-form text strings (e.g. "26.418.24",
#   "  -0.50%  ") rather than numbers, so we force the cell format to Text before
#   assigning -- otherwise Excel auto-parses strings like "1.000" or "243.54" into
#   numeric values and drops the significant trailing zero / formatting.
# - Rows 39/40 (RenderToken / MXToken) additionally swapped ranking order, so their
#   Coin name (B) and Link (C) columns are updated too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    # Force the cell to Text format so numeric-looking strings (prices like
    # "1.000" / "26.418.24") are kept verbatim instead of being parsed as numbers.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    # Reset to the default style so we do not leave a stray Text-format style
    # applied to the cell (the source cells carry no explicit style).
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "26.418.24"
Set-TextCell "E2" "  -0.50%  "

# Row 3
Set-TextCell "D3" "1.724.53"
Set-TextCell "E3" "  -0.25%  "

# Row 4
Set-TextCell "D4" "1.000"
Set-TextCell "E4" "  +0.00%  "

# Row 5
Set-TextCell "D5" "243.54"
Set-TextCell "E5" "  -0.58%  "

# Row 6
Set-TextCell "D6" "1.000"
Set-TextCell "E6" "  -0.01%  "

# Row 7
Set-TextCell "D7" "0.4932"
Set-TextCell "E7" "  +2.60%  "

# Row 8
Set-TextCell "D8" "0.2618"
Set-TextCell "E8" "  -1.85%  "

# Row 9
Set-TextCell "D9" "0.06199"
Set-TextCell "E9" "  +0.33%  "

# Row 10
Set-TextCell "D10" "1.728.14"
Set-TextCell "E10" "  -0.03%  "

# Row 11
Set-TextCell "D11" "0.06998"
Set-TextCell "E11" "  -2.53%  "

# Row 12
Set-TextCell "D12" "15.44"
Set-TextCell "E12" "  -0.80%  "

# Row 13
Set-TextCell "D13" "4.546"
Set-TextCell "E13" "  +0.22%  "

# Row 14
Set-TextCell "D14" "0.5997"
Set-TextCell "E14" "  -1.56%  "

# Row 15
Set-TextCell "D15" "77.48"
Set-TextCell "E15" "  +0.41%  "

# Row 16
Set-TextCell "D16" "0.9999"
Set-TextCell "E16" "  -0.06%  "

# Row 17
Set-TextCell "D17" "26.426.48"
Set-TextCell "E17" "  -0.47%  "

# Row 18
Set-TextCell "D18" "1.000"
Set-TextCell "E18" "  +0.02%  "

# Row 19
Set-TextCell "D19" "0.000007193"
Set-TextCell "E19" "  +3.34%  "

# Row 20
Set-TextCell "E20" "  -1.55%  "

# Row 21
Set-TextCell "D21" "1.949.33"
Set-TextCell "E21" "  -0.24%  "

# Row 22
Set-TextCell "D22" "4.474"
Set-TextCell "E22" "  -0.99%  "

# Row 23
Set-TextCell "D23" "8.583"
Set-TextCell "E23" "  -2.38%  "

# Row 24
Set-TextCell "D24" "5.166"
Set-TextCell "E24" "  -1.28%  "

# Row 25
Set-TextCell "D25" "137.76"
Set-TextCell "E25" "  +0.39%  "

# Row 26
Set-TextCell "E26" "  -0.45%  "

# Row 27
Set-TextCell "E27" "  -0.59%  "

# Row 28
Set-TextCell "D28" "106.94"
Set-TextCell "E28" "  -0.65%  "

# Row 29
Set-TextCell "D29" "1.721"
Set-TextCell "E29" "  -3.08%  "

# Row 30
Set-TextCell "D30" "3.944"
Set-TextCell "E30" "  -0.54%  "

# Row 31
Set-TextCell "D31" "0.08011"
Set-TextCell "E31" "  +0.07%  "

# Row 32
Set-TextCell "D32" "3.679"
Set-TextCell "E32" "  -0.24%  "

# Row 33
Set-TextCell "D33" "0.04524"
Set-TextCell "E33" "  +0.35%  "

# Row 34
Set-TextCell "D34" "0.9993"
Set-TextCell "E34" "  -0.07%  "

# Row 35
Set-TextCell "E35" "  -0.40%  "

# Row 36
Set-TextCell "D36" "0.9986"
Set-TextCell "E36" "  -0.23%  "

# Row 37
Set-TextCell "D37" "0.6273"
Set-TextCell "E37" "  -0.58%  "

# Row 38
Set-TextCell "D38" "0.9517"
Set-TextCell "E38" "  +5.42%  "

# Row 39
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D39" "2.391"
Set-TextCell "E39" "  -0.17%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D40" "1.949"
Set-TextCell "E40" "  -4.65%  "

# Row 41
Set-TextCell "D41" "1.0000"
Set-TextCell "E41" "  -0.15%  "

# Row 42
Set-TextCell "D42" "0.01486"
Set-TextCell "E42" "  -1.14%  "

# Row 43
Set-TextCell "E43" "  -3.57%  "

# Row 44
Set-TextCell "D44" "5.334"
Set-TextCell "E44" "  -2.77%  "

# Row 45
Set-TextCell "D45" "0.3853"
Set-TextCell "E45" "  -0.88%  "

# Row 46
Set-TextCell "D46" "6.819"
Set-TextCell "E46" "  -3.31%  "

# Row 47
Set-TextCell "D47" "0.1167"
Set-TextCell "E47" "  -1.40%  "

# Row 48
Set-TextCell "D48" "0.05371"
Set-TextCell "E48" "  -0.26%  "

# Row 49
Set-TextCell "D49" "7.766"
Set-TextCell "E49" "  -1.00%  "

# Row 50
Set-TextCell "D50" "30.22"
Set-TextCell "E50" "  -1.40%  "

# Row 51
Set-TextCell "E51" "  -1.57%  "
